# TC08_Trials_Filter_Diagnosis-CarcinoUterus.xlsx
# - Add the Neo4j MATCH query text into cell A2 (wraps onto several lines,
#   hence the taller row 2)
# - Update row 2 height to fit the wrapped text
# - Update the sheet's current selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$query = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Carcinosarcoma of the uterus'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

$ws.Range("A2").Value = $query

# Row 2 grows tall enough to show the wrapped query text
$ws.Rows.Item(2).RowHeight = 87

# Move the selection from C6 to B3:B8
$ws.Range("B3:B8").Select()
